# Refresh the cryptos price/volume table (GitHub Actions data pull).
# Note: a handful of "Price" cells hold numeric-looking text (e.g. "23.60",
# "1.00") that must stay text so trailing zeros survive; those are written
# with a leading apostrophe, exactly as typing '23.60 into a cell in Excel
# forces text entry without altering any other formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.352.46'
$ws.Range('E2').Value = '  +3.11%  '
$ws.Range('D3').Value = '2.322.31'
$ws.Range('E3').Value = '  +0.90%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''545.11'
$ws.Range('E5').Value = '  +1.40%  '
$ws.Range('D6').Value = '''130.85'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -1.46%  '
$ws.Range('D9').Value = '2.319.05'
$ws.Range('E9').Value = '  +0.80%  '
$ws.Range('E10').Value = '  +0.41%  '
$ws.Range('E11').Value = '  +0.58%  '
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('E13').Value = '  +0.32%  '
$ws.Range('D14').Value = '''23.60'
$ws.Range('E14').Value = '  -0.80%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '60.329.60'
$ws.Range('E15').Value = '  +3.21%  '
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Value = '2.735.36'
$ws.Range('E16').Value = '  +0.82%  '
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('D18').Value = '2.329.43'
$ws.Range('E18').Value = '  +1.05%  '
$ws.Range('D19').Value = '''10.57'
$ws.Range('E19').Value = '  -0.12%  '
$ws.Range('E20').Value = '  -1.77%  '
$ws.Range('D21').Value = '''314.01'
$ws.Range('E21').Value = '  -0.74%  '
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').Value = '''0.998'
$ws.Range('E23').Value = '  -0.27%  '
$ws.Range('D24').Value = '''63.81'
$ws.Range('E24').Value = '  +1.22%  '
$ws.Range('E25').Value = '  +2.66%  '
$ws.Range('D26').Value = '''1.00'
$ws.Range('E26').Value = '  +0.20%  '
$ws.Range('D27').Value = '''7.85'
$ws.Range('E27').Value = '  -1.31%  '
$ws.Range('E28').Value = '  +3.89%  '
$ws.Range('B29').Value = 'SuiNetwork'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D29').Value = '''1.19'
$ws.Range('E29').Value = '  +9.03%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = '''173.27'
$ws.Range('E30').Value = '  +1.30%  '
$ws.Range('D31').Value = '''1.73'
$ws.Range('E31').Value = '  +0.99%  '
$ws.Range('D32').Value = '0.0₃0732'
$ws.Range('E32').Value = '  +0.83%  '
$ws.Range('E33').Value = '  +1.89%  '
$ws.Range('D34').Value = '''1.37'
$ws.Range('E34').Value = '  +10.31%  '
$ws.Range('E35').Value = '  -0.89%  '
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('D37').Value = '''17.81'
$ws.Range('E37').Value = '  -0.47%  '
$ws.Range('D38').Value = '''0.999'
$ws.Range('E38').Value = '  -0.19%  '
$ws.Range('E39').Value = '  +1.33%  '
$ws.Range('D40').Value = '''323.14'
$ws.Range('E40').Value = '  +10.83%  '
$ws.Range('D41').Value = '''1.54'
$ws.Range('E41').Value = '  +1.93%  '
$ws.Range('D42').Value = '''37.93'
$ws.Range('E42').Value = '  -1.15%  '
$ws.Range('D43').Value = '''137.96'
$ws.Range('E43').Value = '  -2.28%  '
$ws.Range('E44').Value = '  +0.74%  '
$ws.Range('E45').Value = '  -1.34%  '
$ws.Range('D46').Value = '''19.06'
$ws.Range('E46').Value = '  +4.00%  '
$ws.Range('E47').Value = '  +1.01%  '
$ws.Range('E48').Value = '  -0.50%  '
$ws.Range('E49').Value = '  +1.18%  '
$ws.Range('D50').Value = '0.0₆0213'
$ws.Range('E50').Value = '  +17.03%  '
$ws.Range('D51').Value = '''11.02'
$ws.Range('E51').Value = '  +0.59%  '
